$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must remain plain text (the source sheet stores
# every data cell as a string, including numeric-looking price fields like
# "10.60" or "0.592" -- without this, Excel's normal General-format type
# inference would silently turn them into numbers and drop trailing zeros).
function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "58.907.67"
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("D3").Value = "2.574.39"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  -0.01%  "
Set-TextValue "D5" "564.03"
$ws.Range("E5").Value = "  +3.35%  "
Set-TextValue "D6" "142.42"
$ws.Range("E6").Value = "  -0.95%  "
$ws.Range("E7").Value = "  -0.02%  "
Set-TextValue "D8" "0.592"
$ws.Range("E8").Value = "  +1.39%  "
$ws.Range("D9").Value = "2.578.23"
$ws.Range("E9").Value = "  -0.13%  "
Set-TextValue "D10" "6.62"
Set-TextValue "D11" "0.103"
$ws.Range("E11").Value = "  +3.07%  "
Set-TextValue "D12" "0.150"
$ws.Range("E12").Value = "  +7.11%  "
Set-TextValue "D13" "0.341"
$ws.Range("E13").Value = "  +2.92%  "
$ws.Range("D14").Value = "3.025.35"
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").Value = "58.972.26"
$ws.Range("E15").Value = "  +0.97%  "
Set-TextValue "D16" "21.98"
$ws.Range("E16").Value = "  +6.94%  "
Set-TextValue "D17" "0.0000137"
$ws.Range("E17").Value = "  +4.51%  "
$ws.Range("D18").Value = "2.563.87"
$ws.Range("E18").Value = "  -0.53%  "
Set-TextValue "D19" "4.50"
$ws.Range("E19").Value = "  +1.32%  "
Set-TextValue "D20" "334.84"
$ws.Range("E20").Value = "  +0.32%  "
Set-TextValue "D21" "10.14"
$ws.Range("E21").Value = "  +1.33%  "
Set-TextValue "D22" "6.16"
$ws.Range("E22").Value = "  +1.59%  "
Set-TextValue "D23" "0.999"
$ws.Range("E23").Value = "  +0.02%  "
Set-TextValue "D24" "63.65"
$ws.Range("E24").Value = "  -4.28%  "
$ws.Range("E25").Value = "  +5.25%  "
Set-TextValue "D26" "0.999"
$ws.Range("E26").Value = "  +0.06%  "
Set-TextValue "D27" "0.162"
$ws.Range("E27").Value = "  +2.50%  "
Set-TextValue "D28" "7.21"
$ws.Range("E28").Value = "  +2.22%  "
$ws.Range("D29").Value = "0.0₃0777"
$ws.Range("E29").Value = "  +5.60%  "
Set-TextValue "D30" "0.999"
$ws.Range("E30").Value = "  +0.02%  "
Set-TextValue "D31" "1.67"
$ws.Range("E31").Value = "  +1.47%  "
Set-TextValue "D32" "159.64"
$ws.Range("E32").Value = "  +3.23%  "
Set-TextValue "D33" "6.02"
$ws.Range("E33").Value = "  +1.78%  "
Set-TextValue "D34" "18.89"
$ws.Range("E34").Value = "  +0.28%  "
$ws.Range("E35").Value = "  +2.57%  "
Set-TextValue "D36" "0.876"
$ws.Range("E36").Value = "  +2.63%  "
Set-TextValue "D37" "0.877"
$ws.Range("E37").Value = "  +7.62%  "
Set-TextValue "D38" "1.13"
$ws.Range("E38").Value = "  +3.23%  "
Set-TextValue "D39" "36.84"
$ws.Range("E39").Value = "  -0.73%  "
Set-TextValue "D40" "1.48"
$ws.Range("E40").Value = "  +4.27%  "
Set-TextValue "D41" "289.74"
$ws.Range("E41").Value = "  +3.81%  "
Set-TextValue "D42" "3.62"
$ws.Range("E42").Value = "  +1.51%  "
Set-TextValue "D43" "0.998"
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("E44").Value = "  +2.99%  "
Set-TextValue "D45" "0.595"
$ws.Range("E45").Value = "  +0.53%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D46" "0.0533"
$ws.Range("E46").Value = "  +1.07%  "
$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D47" "10.60"
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D48" "19.00"
$ws.Range("E48").Value = "  +2.74%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D49" "125.01"
$ws.Range("E49").Value = "  +12.07%  "
$ws.Range("E50").Value = "  +1.59%  "
Set-TextValue "D51" "18.45"
$ws.Range("E51").Value = "  +4.21%  "
